$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column for existing data rows 2..185 from 45190 -> 45192
for ($r = 2; $r -le 185; $r++) {
    $ws.Cells.Item($r, 3).Value = 45192
}

# Row 185 gains an explicit row height (matches the authored xlsx row attrs)
$ws.Rows.Item(185).RowHeight = 15

# Append the new record as row 186
$ws.Cells.Item(186, 1).Value = "A 44669-2023"
$ws.Cells.Item(186, 2).Value = 45189
$ws.Cells.Item(186, 3).Value = 45192
$ws.Range("B186:C186").NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(186, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(186, 5).Value = "VÅRGÅRDA"
$ws.Cells.Item(186, 7).Value = 0.7
$ws.Cells.Item(186, 8).Value = 0
$ws.Cells.Item(186, 9).Value = 0
$ws.Cells.Item(186, 10).Value = 0
$ws.Cells.Item(186, 11).Value = 0
$ws.Cells.Item(186, 12).Value = 0
$ws.Cells.Item(186, 13).Value = 0
$ws.Cells.Item(186, 14).Value = 0
$ws.Cells.Item(186, 15).Value = 0
$ws.Cells.Item(186, 16).Value = 0
$ws.Cells.Item(186, 17).Value = 0
$ws.Range("R186").WrapText = $true
